$wb = $excel.ActiveWorkbook
$generic = $wb.Worksheets.Item("GENERIC")

# Add the new "mutation" sheet right after "GENERIC"
$ws = $wb.Worksheets.Add($null, $generic)
$ws.Name = "mutation"

# Row 1 - headers
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "name/str"
$ws.Range("C1").Value = "points"
$ws.Range("D1").Value = "description"
$ws.Range("E1").Value = "starting_trait"
$ws.Range("F1").Value = "valid"
$ws.Range("G1").Value = "vitamins_absorb_multi"
$ws.Range("J1").Value = "cancels:list"

# Row 2
$ws.Range("A2").Value = "VEGETARIAN"
$ws.Range("B2").Value = "Meat Intorlerance"
$ws.Range("C2").Value = -2
$ws.Range("D2").Value = "You have problems with eating meat.  It's possible for you to eat it, but you will suffer morale penalties and obtain less nutrition from it."
$ws.Range("E2").Value = $true
$ws.Range("F2").Value = $false
$ws.Range("G2").Value = "flesh"
$ws.Range("J2").Value = "CANNIBAL"

# Row 3
$ws.Range("H3").Value = "vitA"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = "MEATARIAN"

# Row 4
$ws.Range("H4").Value = "vitB"
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = "ANTIFRUIT"

# Row 5
$ws.Range("H5").Value = "vitC"
$ws.Range("I5").Value = 0

# Row 6
$ws.Range("H6").Value = "calcium"
$ws.Range("I6").Value = 0

# Row 7
$ws.Range("H7").Value = "iron"
$ws.Range("I7").Value = 0

# Selection lands on A3 for the new sheet (making it the active/selected tab)
$ws.Range("A3").Select()
